$wb = $excel.ActiveWorkbook

# "Repayment schedule" is the worksheet that holds the variable-instalment
# schedule (3rd tab). A new (blank) column needs to be inserted just before
# the existing "Late" / heading / "Outstanding" columns (column N), pushing
# those three columns one place to the right (N->O, O->P, P->Q).
$ws = $wb.Worksheets.Item("Repayment schedule")

$ws.Columns("N").Insert()

# Inserting a column copies formatting from the column to its left (M); give
# the new column the same width as column M so it matches that behaviour.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and leave the selection on
# cell J16, matching the saved view state.
$ws.Activate() | Out-Null
$ws.Range("J16").Select() | Out-Null
